$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Summary table (rows 2-9) ---
# D3: Estimated SP Sprint 2 for "Lider de desarrollo" 48 -> 42
$ws.Range("D3").Value = 42
# E3: Estimated SP Sprint 3 for "Lider de desarrollo" 8 -> 14
$ws.Range("E3").Value = 14

# --- Detail table (rows 20-27) ---
# E21: Estimated SP Sprint 2 for "Lider de desarrollo" 48 -> 42 (mirrors D3)
$ws.Range("E21").Value = 42
# G21: Estimated SP Sprint 3 for "Lider de desarrollo" 8 -> 14 (mirrors E3)
$ws.Range("G21").Value = 14

# Column F = Real SP Sprint 2
$ws.Range("F20").Value = 44
$ws.Range("F21").Value = 61
$ws.Range("F22").Value = 60
$ws.Range("F23").Value = 90
$ws.Range("F24").Value = 35
$ws.Range("F25").Value = 70
$ws.Range("F26").Value = 30
$ws.Range("F27").Formula = "=SUM(F20:F26)"

# Update selection to match author's final cursor position
$ws.Range("F28").Select()

$wb.Save()
